# master_story sheet cleanup:
#  - drop the "表 1" title row (was row 1, merged B1:F1)
#  - drop the now-empty leading index column (was column A)
# This shifts the header row (type/msg_data/arg1/arg2/arg3) from B2:F2
# up to A1:E1, and the numbered message rows up/left so the message
# text lands in column B instead of C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Delete()
$ws.Columns.Item(1).Delete()

# the old layout carried blank styled filler cells out to column F;
# clear those leftovers now that the real data only spans A:B
$ws.Range("C2:E18").Clear()

# re-tune the print setup for the now-narrower sheet: plain A4 @ 100%
# instead of the old "squeeze everything onto one page" fit-to-page hack
$ps = $ws.PageSetup
$ps.LeftMargin = 0.7 * 72
$ps.RightMargin = 0.7 * 72
$ps.TopMargin = 0.75 * 72
$ps.BottomMargin = 0.75 * 72
$ps.HeaderMargin = 0.3 * 72
$ps.FooterMargin = 0.3 * 72
$ps.PaperSize = 9
$ps.Orientation = 1
$ps.Zoom = 100
